$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.855.07"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "1.808.24"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.79"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "38.86"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -8.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.318"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.48%  "
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "2.070.00"
$ws.Range("E12").Value = "  -3.01%  "
$ws.Range("D13").Value = "1.820.05"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.659"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -3.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "10.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -7.17%  "
$ws.Range("E16").Value = "  -4.44%  "
$ws.Range("D17").Value = "34.841.85"
$ws.Range("E17").Value = "  -2.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.21"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").Value = "0.0₃0779"
$ws.Range("E19").Value = "  -3.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.79"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -5.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.62"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.94%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "173.56"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.14"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.49%  "
$ws.Range("E28").Value = "  -3.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.50"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.79%  "
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.97"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0544"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.92"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.17"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.85%  "
$ws.Range("E35").Value = "  -8.08%  "
$ws.Range("E36").Value = "  -1.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "90.77"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -11.46%  "
$ws.Range("E38").Value = "  +5.90%  "
$ws.Range("D39").Value = "1.308.11"
$ws.Range("E39").Value = "  -4.47%  "
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.953"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "14.10"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -12.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.69"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("E46").Value = "  -2.18%  "
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("D48").Value = "1.998.18"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  +7.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.60"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.11%  "
